# Applies the changes described by the commit:
#   1. Slide 5's table switches to the built-in "Medium Style 2 - Accent 1"
#      table style ({2F136414-0740-4E08-A5F2-18B06836D584}).
#   2. The presentation's theme (ppt/theme/theme1.xml, driving the slide
#      master / slides) is swapped from the custom "Integral" (Red Violet)
#      colour scheme to the stock "Office Theme" colour scheme.
#
# Note: the notes-master theme part (ppt/theme/theme2.xml) is not
# independently addressable through the PowerPoint COM object model
# exposed here (NotesMaster.Theme / Master.Theme all resolve back to the
# single presentation theme object), so only the reachable half of the
# swap (the slide/master-facing theme) can be reproduced via COM.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------

$s = $p.Slides.Item(5)
$tableShape = $s.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{2F136414-0740-4E08-A5F2-18B06836D584}")

# --- 2. Swap theme colour scheme: Integral (Red Violet) -> Office Theme --

$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# dk1 = 000000
$cs.Item(1).RGB = 0x00 + (0x00 * 256) + (0x00 * 65536)
# lt1 = FFFFFF
$cs.Item(2).RGB = 0xFF + (0xFF * 256) + (0xFF * 65536)
# dk2 = 44546A
$cs.Item(3).RGB = 0x44 + (0x54 * 256) + (0x6A * 65536)
# lt2 = E7E6E6
$cs.Item(4).RGB = 0xE7 + (0xE6 * 256) + (0xE6 * 65536)
# accent1 = 5B9BD5
$cs.Item(5).RGB = 0x5B + (0x9B * 256) + (0xD5 * 65536)
# accent2 = ED7D31
$cs.Item(6).RGB = 0xED + (0x7D * 256) + (0x31 * 65536)
# accent3 = A5A5A5
$cs.Item(7).RGB = 0xA5 + (0xA5 * 256) + (0xA5 * 65536)
# accent4 = FFC000
$cs.Item(8).RGB = 0xFF + (0xC0 * 256) + (0x00 * 65536)
# accent5 = 4472C4
$cs.Item(9).RGB = 0x44 + (0x72 * 256) + (0xC4 * 65536)
# accent6 = 70AD47
$cs.Item(10).RGB = 0x70 + (0xAD * 256) + (0x47 * 65536)
# hlink = 0563C1
$cs.Item(11).RGB = 0x05 + (0x63 * 256) + (0xC1 * 65536)
# folHlink = 954F72
$cs.Item(12).RGB = 0x95 + (0x4F * 256) + (0x72 * 65536)

Write-Output "done"
